$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Final target data for rows 2-17 (columns A=Player, B=Position, C=Team)
$data = @(
    @("Tyler Herro",       "PG,SG",    "Miami Heat"),
    @("Klay Thompson",     "SG,SF",    "Dallas Mavericks"),
    @("Mikal Bridges",     "SG,SF,PF", "New York Knicks"),
    @("De'Aaron Fox",      "PG",       "Sacramento Kings"),
    @("Jonathan Kuminga",  "SF,PF",    "Golden State Warriors"),
    @("DeMar DeRozan",     "SF,PF",    "Sacramento Kings"),
    @("Nikola Vucevic",    "PF,C",     "Chicago Bulls"),
    @("Santi Aldama",      "PF,C",     "Memphis Grizzlies"),
    @("Evan Mobley",       "PF,C",     "Cleveland Cavaliers"),
    @("Brook Lopez",       "C",        "Milwaukee Bucks"),
    @("Josh Giddey",       "PG,SG,SF", "Chicago Bulls"),
    @("P.J. Washington",   "PF",       "Dallas Mavericks"),
    @("Ja Morant",         "PG",       "Memphis Grizzlies"),
    @("Miles Bridges",     "SF,PF",    "Charlotte Hornets"),
    @("Scottie Barnes",    "SG,SF,PF", "Toronto Raptors"),
    @("Luka Doncic",       "PG,SG",    "Dallas Mavericks")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
